$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.167.51'
$ws.Range("E2").Value = '  -1.37%  '
$ws.Range("D3").Value = '2.229.62'
$ws.Range("E3").Value = '  -2.24%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.634'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.86%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '62.58'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.55%  '
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.436'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0948'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -8.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.62'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.83%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.39'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("E13").Value = '  +0.23%  '
$ws.Range("D14").Value = '2.568.24'
$ws.Range("E14").Value = '  -1.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.25'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.96'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.818'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.77%  '
$ws.Range("D18").Value = '2.241.69'
$ws.Range("E18").Value = '  -1.55%  '
$ws.Range("D19").Value = '43.176.76'
$ws.Range("E19").Value = '  -1.23%  '
$ws.Range("D20").Value = '0.0₃0953'
$ws.Range("E20").Value = '  -4.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '243.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.78%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.67'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +31.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.39'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.99%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.65'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '171.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.21'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.130'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.38'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.124'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.81'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0664'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.81'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.57'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.27'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.23'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0247'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.00%  '
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.48'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.45'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.03%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.82'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.98%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '96.21'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.51%  '
$ws.Range("E46").Value = '  -3.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0925'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.61%  '
$ws.Range("D48").Value = '1.439.01'
$ws.Range("E48").Value = '  -2.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000206'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.75'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.25'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.26%  '
